$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same style as other header cells (copy format from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:50:18.712997"
$ws.Range("F3").Value = "2021-10-05 10:50:18.713008"
$ws.Range("F4").Value = "2021-10-05 10:50:18.713012"
